$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Coin, Link, Price, Volume(1h)) per row.
# A leading apostrophe forces Excel to store a value as text (quote-prefix)
# instead of auto-converting numeric-looking strings like "1.00" into numbers;
# Excel strips the apostrophe itself and keeps the literal text as the value.
$rows = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '67.746.19'; E = '  +0.17%  ' }
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '3.791.82'; E = '  +0.37%  ' }
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '''0.998'; E = '  -0.25%  ' }
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '''599.02'; E = '  +0.69%  ' }
    @{ Row = 6; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '''165.15'; E = '  -1.01%  ' }
    @{ Row = 7; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '''1.00'; E = '  +0.03%  ' }
    @{ Row = 8; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '''0.517'; E = '  -0.43%  ' }
    @{ Row = 9; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '''0.159'; E = '  -0.05%  ' }
    @{ Row = 10; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '''0.453'; E = '  +0.97%  ' }
    @{ Row = 11; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '''6.45'; E = '  +2.38%  ' }
    @{ Row = 12; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '''0.0000249'; E = '  -1.80%  ' }
    @{ Row = 13; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '''35.77'; E = '  -0.69%  ' }
    @{ Row = 14; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '4.429.55'; E = '  +0.46%  ' }
    @{ Row = 15; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '3.813.43'; E = '  +0.83%  ' }
    @{ Row = 16; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '67.798.06'; E = '  +0.28%  ' }
    @{ Row = 17; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '''18.40'; E = '  -0.53%  ' }
    @{ Row = 18; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '''0.113'; E = '  +1.56%  ' }
    @{ Row = 19; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '''7.06'; E = '  +0.69%  ' }
    @{ Row = 20; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '''463.29'; E = '  +0.90%  ' }
    @{ Row = 21; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '''9.80'; E = '  -2.35%  ' }
    @{ Row = 22; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '''0.701'; E = '  +0.55%  ' }
    @{ Row = 23; B = 'PEPE'; C = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D = '''0.0000147'; E = '  -5.20%  ' }
    @{ Row = 24; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '''82.84'; E = '  -0.49%  ' }
    @{ Row = 25; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '''12.02'; E = '  +0.15%  ' }
    @{ Row = 26; B = 'Fetch.AI'; C = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D = '''2.10'; E = '  -0.53%  ' }
    @{ Row = 27; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '''10.02'; E = '  -0.13%  ' }
    @{ Row = 28; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '''1.00'; E = '  -0.07%  ' }
    @{ Row = 29; B = 'WrappedeETH'; C = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; D = '3.941.17'; E = '  +0.41%  ' }
    @{ Row = 30; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '''2.76'; E = '  -0.26%  ' }
    @{ Row = 31; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '''7.41'; E = '  +2.61%  ' }
    @{ Row = 32; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '''2.21'; E = '  -1.02%  ' }
    @{ Row = 33; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '''29.24'; E = '  -1.13%  ' }
    @{ Row = 34; B = 'Binance-PegBSC-USD'; C = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D = '''1.00'; E = '  -0.03%  ' }
    @{ Row = 35; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '''9.04'; E = '  -0.50%  ' }
    @{ Row = 36; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '''0.0994'; E = '  -0.72%  ' }
    @{ Row = 37; B = 'Kaspa'; C = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D = '''0.139'; E = '  +1.01%  ' }
    @{ Row = 38; B = 'dogwifhat'; C = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D = '''3.24'; E = '  -3.90%  ' }
    @{ Row = 39; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '''5.76'; E = '  -0.14%  ' }
    @{ Row = 40; B = 'Mantle'; C = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D = '''0.985'; E = '  -0.67%  ' }
    @{ Row = 41; B = 'FirstDigitalUSD'; C = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D = '''1.00'; E = '  +0.03%  ' }
    @{ Row = 42; B = 'USDe'; C = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; D = '''1.00'; E = '  +0.01%  ' }
    @{ Row = 43; B = 'Arweave'; C = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'; D = '''44.82'; E = '  -1.52%  ' }
    @{ Row = 44; B = 'OKB'; C = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D = '''47.67'; E = '  -0.78%  ' }
    @{ Row = 45; B = 'TheGraph'; C = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D = '''0.298'; E = '  +0.25%  ' }
    @{ Row = 46; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '''151.40'; E = '  +1.05%  ' }
    @{ Row = 47; B = 'ONDO'; C = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; D = '''1.39'; E = '  +9.11%  ' }
    @{ Row = 48; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '''8.34'; E = '  +0.65%  ' }
    @{ Row = 49; B = 'Bittensor'; C = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D = '''397.60'; E = '  +1.08%  ' }
    @{ Row = 50; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '''27.18'; E = '  +1.42%  ' }
    @{ Row = 51; B = 'Stacks'; C = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D = '''1.85'; E = '  +1.78%  ' }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

